# Update the "© 2015 The MathWorks, Inc." copyright notice on the Title
# Slide layout to reference 2016 instead, per:
#   "updated PowerPoint copyright to 2016"
#
# The Copyright text box is a userDrawn shape that lives on the "Title
# Slide" custom layout (it is not a per-slide shape / placeholder), so we
# reach it through the slide master's CustomLayouts collection rather
# than through any particular Slide's Shapes collection.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$titleLayout = $master.CustomLayouts.Item(1)
$copyright = $titleLayout.Shapes.Item("Copyright")
$copyrightRange = $copyright.TextFrame.TextRange

$copyrightRange.Text = "© 2016 The MathWorks, Inc."
